{"js": "const oldText = \"except the power to consent to guardianship, adoption, or marriage.\";\nconst newText = \"except the power to consent to adoption or to release for adoption.\";\n\nconst searchResults = context.document.body.search(oldText, { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found: \" + oldText);\n}\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"except the power to consent to guardianship, adoption, or marriage.\"\n$newText = \"except the power to consent to adoption or to release for adoption.\"\n\n# Locate the paragraph containing the sentence to be updated (paragraph 3 of\n# the template body, per the commit message) and replace just that sentence,\n# scoping the Find/Replace to the paragraph's own range so no other part of\n# the document is touched.\n$found = $false\nforeach ($p in $d.Paragraphs) {\n    $pr = $p.Range\n    if ($pr.Text -like \"*$oldText*\") {\n        $f = $pr.Find\n        $f.ClearFormatting()\n        $f.Text = $oldText\n        $f.Replacement.ClearFormatting()\n        $f.Replacement.Text = $newText\n        $f.Execute([ref]$f.Text, $true, $true, $false, $false, $false, $true, 1, $false, $f.Replacement.Text, 2) | Out-Null\n        $found = $true\n        break\n    }\n}\n\nif (-not $found) {\n    throw \"Target sentence not found: $oldText\"\n}\n"}
